$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Swap the "Periodo Mora" (E) and "Valor Mora" (F) values between row 16 and row 17
$ws.Range("E16").Value = "2304"
$ws.Range("F16").Value = 41822
$ws.Range("E17").Value = "2303"
$ws.Range("F17").Value = 52000
